$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 369; everything below shifts down by one.
$ws.Rows.Item(369).EntireRow.Insert()

# Populate the newly inserted row 369 with the new data record
$ws.Cells.Item(369, 1).Value = 8
$ws.Cells.Item(369, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(369, 3).Value = "Coquimbo"
$ws.Cells.Item(369, 4).Value = 45173
$ws.Cells.Item(369, 5).Value = 4
$ws.Cells.Item(369, 6).Value = 100112031
$ws.Cells.Item(369, 7).Value = "Poroto verde"
$ws.Cells.Item(369, 8).Value = "Magnum"
$ws.Cells.Item(369, 9).Value = "Primera"
$ws.Cells.Item(369, 10).Value = 400
$ws.Cells.Item(369, 11).Value = 27000
$ws.Cells.Item(369, 12).Value = 28000
$ws.Cells.Item(369, 13).Value = 27500
$ws.Cells.Item(369, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(369, 15).Value = "Perú"
$ws.Cells.Item(369, 16).Value = 1100
$ws.Cells.Item(369, 17).Value = 25
$ws.Cells.Item(369, 18).Value = "Hortaliza"
